$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.754.95"
$ws.Range("E2").Value = "  +2.24%  "
$ws.Range("D3").Value = "3.383.27"
$ws.Range("E3").Value = "  +1.26%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "582.32"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.35%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "179.98"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +1.25%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.596"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.66%  "
$ws.Range("E9").Value = "  +8.86%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.590"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.18%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "48.50"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.61%  "
$ws.Range("E12").Value = "  +4.37%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "683.25"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -1.35%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.64"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +2.25%  "
$ws.Range("D15").Value = "3.925.23"
$ws.Range("E15").Value = "  +1.06%  "
$ws.Range("D16").Value = "69.707.28"
$ws.Range("E16").Value = "  +2.03%  "
$ws.Range("E17").Value = "  +1.05%  "
$ws.Range("D18").Value = "3.382.46"
$ws.Range("E18").Value = "  +1.33%  "
$ws.Range("E19").Value = "  +1.25%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.29"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.89%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.914"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +2.02%  "
$ws.Range("E22").Value = "  +2.17%  "
$ws.Range("E23").Value = "  -1.59%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "101.89"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.58%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.90"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.56%  "
$ws.Range("E26").Value = "  +0.15%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.75"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +2.49%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "33.65"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.61%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.76"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +2.53%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.93"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.94%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.86"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +15.28%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.08"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.26%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "556.87"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -1.65%  "
$ws.Range("E34").Value = "  +0.62%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "57.95"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.81%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("D37").Value = "3.607.26"
$ws.Range("E37").Value = "  -2.78%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.140"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +2.67%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "35.39"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.37%  "
$ws.Range("D40").Value = "0.0₃0730"
$ws.Range("E40").Value = "  +8.51%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.76"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +5.19%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.32"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +3.93%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0427"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +3.45%  "
$ws.Range("E44").Value = "  +0.39%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.67"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.34%  "
$ws.Range("E46").Value = "  +0.16%  "
$ws.Range("E47").Value = "  +0.03%  "
$ws.Range("E48").Value = "  +3.89%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "130.36"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.17%  "
$ws.Range("E50").Value = "  +1.25%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.49"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.61%  "
